$wb = $excel.ActiveWorkbook

# Add the new "rho" worksheet as the last tab (after "Adjusted2")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "rho"

# Copy the header formatting (bold, centered, thin-bordered) from an existing
# sheet's header row so the new sheet reuses the same cell style.
$srcHeader = $wb.Worksheets.Item("Adjusted").Range("A1:F1")
$srcHeader.Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# Header row text
$headers = @("theta", "Jxx", "Jyy", "beta", "gamma", "trace_sq")
for ($c = 1; $c -le $headers.Count; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Data rows: theta, Jxx, Jyy, beta, gamma, trace_sq
$data = @(
    @(0,0.00193466609378859,0.9980653681651108,-0.0159222707822376,0,0.9966452594771784),
    @(5,0.03774810522500854,0.9622518396249239,-0.1905857933890676,0,0.9999994115931092),
    @(10,0.1317658986024291,0.8682341433791049,-0.3347319383646786,0,0.9952837208865037),
    @(15,0.2698379312570733,0.73016214786142,-0.4356074999044452,0,0.9854570592607016),
    @(20,0.4258182378288504,0.5741817461434889,-0.4835477452990676,0,0.9786426932396772),
    @(25,0.6016459029625263,0.3983541098866221,-0.4835001262464965,0,0.9882085335759128),
    @(30,0.7570303134257383,0.2429696919489803,-0.4237348620356136,0,0.9912316332599351),
    @(35,0.8857858169408038,0.1142329877294931,-0.3180972493350357,0,1.000037409048126),
    @(40,0.9712261945637656,0.02877422233326949,-0.1671713270305639,0,1.000000782040017),
    @(45,0.9980284094936002,0.00197160076849049,0.001043480351031901,0,0.9960667710684015),
    @(50,0.9624685513557465,0.03753156866800373,0.1737318873195663,0,0.9881198683387473),
    @(55,0.8732650845918691,0.1267349424279521,0.3217522916708382,0,0.9857027279903331),
    @(60,0.7320320719302075,0.267968146975384,0.438943222663999,0,0.9930201875731675),
    @(65,0.5788102362977725,0.4211898934090182,0.4918298635613694,0,0.9962154453345737),
    @(70,0.4141676533940798,0.5858325329549383,0.4899534589677298,0,0.9948433856952434),
    @(75,0.2498724551897367,0.75012761737949,0.4276994652928388,0,0.9909813514415378),
    @(80,0.1203202736650861,0.8796797687904249,0.3166513794355037,0,0.9888496560708313),
    @(85,0.03209224622152107,0.9679080212922321,0.1653720777968055,0,0.9925716981790521),
    @(90,0.001897114660365675,0.9981028999871777,-0.01019143049454269,0,0.9964207285178989)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
    }
}

$null = $ws.Range("A1").Select()
